$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new DAMSLTag (col I) / DialogAct (col J)
$updates = @(
    @{ Row = 11; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 23; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 25; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 26; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 33; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 38; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 49; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 52; Tag = "%"; Act = "Uninterpretable" },
    @{ Row = 61; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 75; Tag = "ba"; Act = "Appreciation" },
    @{ Row = 83; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 87; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 88; Tag = "%"; Act = "Uninterpretable" },
    @{ Row = 89; Tag = "sd"; Act = "Statement-non-opinion" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.Tag
    $ws.Cells.Item($u.Row, 10).Value = $u.Act
}

$wb.Save()
